$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").NumberFormat = "@"
$ws.Range("D2").Value = "283.25"
$ws.Range("E2").NumberFormat = "@"
$ws.Range("E2").Value = "1.87%"
$ws.Range("D3").NumberFormat = "@"
$ws.Range("D3").Value = "28.34"
$ws.Range("E3").NumberFormat = "@"
$ws.Range("E3").Value = "3.97%"
$ws.Range("D4").NumberFormat = "@"
$ws.Range("D4").Value = "5.052"
$ws.Range("E4").NumberFormat = "@"
$ws.Range("E4").Value = "3.87%"
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "0.06516"
$ws.Range("E5").NumberFormat = "@"
$ws.Range("E5").Value = "1.46%"
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = "7.230"
$ws.Range("E6").NumberFormat = "@"
$ws.Range("E6").Value = "3.13%"
$ws.Range("B7").NumberFormat = "@"
$ws.Range("B7").Value = "FTXToken"
$ws.Range("C7").NumberFormat = "@"
$ws.Range("C7").Value = "https://coinranking.com/coin/NfeOYfNcl+ftxtoken-ftt"
$ws.Range("D7").NumberFormat = "@"
$ws.Range("D7").Value = "1.386"
$ws.Range("E7").NumberFormat = "@"
$ws.Range("E7").Value = "16.40%"
$ws.Range("B8").NumberFormat = "@"
$ws.Range("B8").Value = "MXToken"
$ws.Range("C8").NumberFormat = "@"
$ws.Range("C8").Value = "https://coinranking.com/coin/QUC5kVAxSoB-+mxtoken-mx"
$ws.Range("D8").NumberFormat = "@"
$ws.Range("D8").Value = "0.9184"
$ws.Range("E8").NumberFormat = "@"
$ws.Range("E8").Value = "3.60%"
$ws.Range("B9").NumberFormat = "@"
$ws.Range("B9").Value = "WazirX"
$ws.Range("C9").NumberFormat = "@"
$ws.Range("C9").Value = "https://coinranking.com/coin/6QK-8hUZ+wazirx-wrx"
$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = "0.1535"
$ws.Range("E9").NumberFormat = "@"
$ws.Range("E9").Value = "-0.28%"
$ws.Range("B10").NumberFormat = "@"
$ws.Range("B10").Value = "LiechtensteinCryptoassetsExchange"
$ws.Range("C10").NumberFormat = "@"
$ws.Range("C10").Value = "https://coinranking.com/coin/v4IW9oaF+liechtensteincryptoassetsexchange-lcx"
$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = "0.06630"
$ws.Range("E10").NumberFormat = "@"
$ws.Range("E10").Value = "28.34%"
$ws.Range("B11").NumberFormat = "@"
$ws.Range("B11").Value = "MandalaExchangeToken"
$ws.Range("C11").NumberFormat = "@"
$ws.Range("C11").Value = "https://coinranking.com/coin/lviNIbma2Xuqs+mandalaexchangetoken-mdx"
$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = "0.07542"
$ws.Range("E11").NumberFormat = "@"
$ws.Range("E11").Value = "0.43%"
$ws.Range("B12").NumberFormat = "@"
$ws.Range("B12").Value = "BitrueCoin"
$ws.Range("C12").NumberFormat = "@"
$ws.Range("C12").Value = "https://coinranking.com/coin/SLYjzF4ty+bitruecoin-btr"
$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = "0.02799"
$ws.Range("E12").NumberFormat = "@"
$ws.Range("E12").Value = "-3.19%"
$ws.Range("B13").NumberFormat = "@"
$ws.Range("B13").Value = "BitMartToken"
$ws.Range("C13").NumberFormat = "@"
$ws.Range("C13").Value = "https://coinranking.com/coin/6uzcPMFgWUJNH+bitmarttoken-bmx"
$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = "0.08975"
$ws.Range("E13").NumberFormat = "@"
$ws.Range("E13").Value = "-0.08%"
$ws.Range("B14").NumberFormat = "@"
$ws.Range("B14").Value = "BitForexToken"
$ws.Range("C14").NumberFormat = "@"
$ws.Range("C14").Value = "https://coinranking.com/coin/2nh5ugplNocUp+bitforextoken-bf"
$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = "0.001585"
$ws.Range("E14").NumberFormat = "@"
$ws.Range("E14").Value = "1.46%"
$ws.Range("B15").NumberFormat = "@"
$ws.Range("B15").Value = "One"
$ws.Range("C15").NumberFormat = "@"
$ws.Range("C15").Value = "https://coinranking.com/coin/6Lga5NiXX3rT+one-one"
$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = "0.0006379"
$ws.Range("E15").NumberFormat = "@"
$ws.Range("E15").Value = "0.26%"
$ws.Range("B16").NumberFormat = "@"
$ws.Range("B16").Value = "TigerCash"
$ws.Range("C16").NumberFormat = "@"
$ws.Range("C16").Value = "https://coinranking.com/coin/6hIn06L2+tigercash-tch"
$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = "0.006150"
$ws.Range("E16").NumberFormat = "@"
$ws.Range("E16").Value = "0.17%"
$ws.Range("B17").NumberFormat = "@"
$ws.Range("B17").Value = "LEO"
$ws.Range("C17").NumberFormat = "@"
$ws.Range("C17").Value = "https://coinranking.com/coin/mqtUpyBxu8O8+leo-leo"
$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = "3.446"
$ws.Range("E17").NumberFormat = "@"
$ws.Range("E17").Value = "-0.75%"
$ws.Range("B18").NumberFormat = "@"
$ws.Range("B18").Value = "GateToken"
$ws.Range("C18").NumberFormat = "@"
$ws.Range("C18").Value = "https://coinranking.com/coin/t7m8DZVyMsAu+gatetoken-gt"
$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = "3.351"
$ws.Range("E18").NumberFormat = "@"
$ws.Range("E18").Value = "1.46%"
$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = "2.237"
$ws.Range("E19").NumberFormat = "@"
$ws.Range("E19").Value = "-1.50%"
$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = "3.986"
$ws.Range("E22").NumberFormat = "@"
$ws.Range("E22").Value = "2.11%"
$ws.Range("E23").NumberFormat = "@"
$ws.Range("E23").Value = "1.63%"
$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = "0.04431"
$ws.Range("E24").NumberFormat = "@"
$ws.Range("E24").Value = "0.65%"
$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = "0.001183"
$ws.Range("E25").NumberFormat = "@"
$ws.Range("E25").Value = "0.63%"
$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = "0.004435"
$ws.Range("E26").NumberFormat = "@"
$ws.Range("E26").Value = "14.12%"
$ws.Range("E27").NumberFormat = "@"
$ws.Range("E27").Value = "1.67%"
$ws.Range("E28").NumberFormat = "@"
$ws.Range("E28").Value = "-1.57%"
$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = "0.04116"
$ws.Range("E40").NumberFormat = "@"
$ws.Range("E40").Value = "-0.09%"
$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = "0.006673"
$ws.Range("E41").NumberFormat = "@"
$ws.Range("E41").Value = "-2.23%"
$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = "0.1229"
$ws.Range("E42").NumberFormat = "@"
$ws.Range("E42").Value = "4.69%"
$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = "0.002149"
$ws.Range("E43").NumberFormat = "@"
$ws.Range("E43").Value = "13.74%"
$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = "0.01208"
$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = "0.00005687"
$ws.Range("E45").NumberFormat = "@"
$ws.Range("E45").Value = "6.67%"
$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = "1.966"
$ws.Range("E46").NumberFormat = "@"
$ws.Range("E46").Value = "39.89%"
